$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.913.88'
$ws.Range("E2").Value = '  -2.95%  '

$ws.Range("D3").Value = '1.676.54'
$ws.Range("E3").Value = '  -3.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3660'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3358'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07292'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.158'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.791'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.64%  '

$ws.Range("D16").Value = '1.671.63'
$ws.Range("E16").Value = '  -3.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001096'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.11%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9982'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06587'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '81.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.156'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").Value = '24.856.22'
$ws.Range("E24").Value = '  -3.14%  '

$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.675'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.256'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.66%  '

$ws.Range("D31").Value = '1.863.38'
$ws.Range("E31").Value = '  -3.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.432'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.144'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08581'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.729'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.404'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06445'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02332'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.666'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2154'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.246'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6242'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9987'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.786'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5949'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.032'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07151'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.35%  '
